$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# --- "2) Change Data" -> "4) " + "Change Data" (two runs) ---
$sh1 = $s.Shapes.Item("TextBox 33")
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(1, 3).Text = "4) "

# --- "3a) Commit Transaction" -> "6a" + ") Commit Transaction" (two runs) ---
$sh2 = $s.Shapes.Item("TextBox 35")
$tr2 = $sh2.TextFrame.TextRange
$tr2.Characters(1, 2).Text = "6a"

# --- "3b) Rollback Transaction" -> "6b" + ") Rollback Transaction" (two runs) ---
$sh3 = $s.Shapes.Item("TextBox 37")
$tr3 = $sh3.TextFrame.TextRange
$tr3.Characters(1, 2).Text = "6b"

# --- "Adding data to CSV-File during add \r and update is not shown." textbox ---
$sh4 = $s.Shapes.Item("TextBox 100")

# Reposition / resize the shape
$sh4.Left = 386.72133858267716
$sh4.Width = 166.0559842519685

$tr4 = $sh4.TextFrame.TextRange

# Paragraph 1: "Adding data to CSV-File during add " -> "Adding data to CSV-File during "
$tr4.Characters(1, 35).Text = "Adding data to CSV-File during "

# Paragraph 2, run 1: "and update i" -> "transaction i" (will later be split into two runs)
$tr4.Characters(33, 12).Text = "transaction i"

# Split the trailing "i" off into its own run
$tr4.Characters(45, 1).Text = "i"
